$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells are treated as text so values like "1.00" or "238.45"
# are not coerced into numbers (which would drop formatting / trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "90.504.93"
$ws.Range("E2").Value = "  +0.93%  "

$ws.Range("D3").Value = "3.145.17"
$ws.Range("E3").Value = "  +2.51%  "

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").Value = "238.45"
$ws.Range("E5").Value = "  +0.60%  "

$ws.Range("D6").Value = "615.75"
$ws.Range("E6").Value = "  -0.21%  "

$ws.Range("E7").Value = "  +4.95%  "

$ws.Range("D8").Value = "0.372"
$ws.Range("E8").Value = "  +2.42%  "

$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").Value = "3.142.23"
$ws.Range("E10").Value = "  +2.46%  "

$ws.Range("D11").Value = "0.738"
$ws.Range("E11").Value = "  +4.42%  "

$ws.Range("E12").Value = "  +1.49%  "

$ws.Range("D13").Value = "0.0000246"
$ws.Range("E13").Value = "  +0.11%  "

$ws.Range("D14").Value = "34.89"
$ws.Range("E14").Value = "  +0.67%  "

$ws.Range("D15").Value = "5.54"
$ws.Range("E15").Value = "  +3.41%  "

$ws.Range("D16").Value = "90.691.22"
$ws.Range("E16").Value = "  +1.53%  "

$ws.Range("D17").Value = "3.729.93"
$ws.Range("E17").Value = "  +2.54%  "

$ws.Range("D18").Value = "3.188.42"
$ws.Range("E18").Value = "  +3.81%  "

$ws.Range("D19").Value = "3.70"
$ws.Range("E19").Value = "  -1.18%  "

$ws.Range("D20").Value = "15.06"
$ws.Range("E20").Value = "  +9.58%  "

$ws.Range("D21").Value = "5.98"
$ws.Range("E21").Value = "  +10.86%  "

$ws.Range("D22").Value = "448.61"
$ws.Range("E22").Value = "  +4.39%  "

$ws.Range("D23").Value = "0.0000202"
$ws.Range("E23").Value = "  -4.66%  "

$ws.Range("D24").Value = "9.07"
$ws.Range("E24").Value = "  +4.58%  "

$ws.Range("D25").Value = "5.94"
$ws.Range("E25").Value = "  +7.05%  "

$ws.Range("D26").Value = "88.49"
$ws.Range("E26").Value = "  +1.98%  "

$ws.Range("D27").Value = "11.89"
$ws.Range("E27").Value = "  +2.13%  "

$ws.Range("D28").Value = "3.328.25"
$ws.Range("E28").Value = "  +2.97%  "

$ws.Range("D29").Value = "1.00"

$ws.Range("D30").Value = "0.139"
$ws.Range("E30").Value = "  +53.99%  "

$ws.Range("D31").Value = "0.234"
$ws.Range("E31").Value = "  +17.41%  "

$ws.Range("D32").Value = "0.170"
$ws.Range("E32").Value = "  +8.46%  "

$ws.Range("D33").Value = "9.26"
$ws.Range("E33").Value = "  +3.63%  "

$ws.Range("D34").Value = "0.173"
$ws.Range("E34").Value = "  +16.01%  "

$ws.Range("E35").Value = "  -5.92%  "

$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "26.30"
$ws.Range("E36").Value = "  +3.01%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").Value = "7.59"
$ws.Range("E37").Value = "  +6.58%  "

$ws.Range("E38").Value = "  +5.17%  "

$ws.Range("D39").Value = "505.99"
$ws.Range("E39").Value = "  +3.12%  "

$ws.Range("D40").Value = "1.33"
$ws.Range("E40").Value = "  +6.35%  "

$ws.Range("D41").Value = "3.83"
$ws.Range("E41").Value = "  -4.28%  "

$ws.Range("D42").Value = "0.444"
$ws.Range("E42").Value = "  +11.73%  "

$ws.Range("D43").Value = "3.43"
$ws.Range("E43").Value = "  -4.93%  "

$ws.Range("D44").Value = "22.09"
$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("D46").Value = "0.713"
$ws.Range("E46").Value = "  +6.18%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "156.50"
$ws.Range("E47").Value = "  +0.03%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "1.92"
$ws.Range("E48").Value = "  +3.81%  "

$ws.Range("E49").Value = "  +5.96%  "

$ws.Range("E50").Value = "  +3.91%  "

$ws.Range("D51").Value = "43.90"
$ws.Range("E51").Value = "  -1.05%  "

# Reset column D style back to the default (no explicit number format)
# now that the values are committed as text, matching the original workbook styling.
$ws.Range("D2:D51").Style = "Normal"
